$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.404345
$ws.Range("H2").Value = 40.213035
$ws.Range("I2").Value = 0.01122005832922476
$ws.Range("J2").Value = 0.01122005832922476
$ws.Range("M2").Value = 1.949849666666667
$ws.Range("N2").Value = 5.849549000000001
$ws.Range("O2").Value = 0.06676506732104066
$ws.Range("P2").Value = 0.06676506732104066
$ws.Range("Q2").Value = 26.136457630135
$ws.Range("R2").Value = 235.228118671215
$ws.Range("S2").Value = 0.0007491079496966941
$ws.Range("T2").Value = 0.0007491079496966941
$ws.Range("G3").Value = 13.404345
$ws.Range("H3").Value = 40.213035
$ws.Range("I3").Value = 0.01122005832922476
$ws.Range("J3").Value = 0.01122005832922476
$ws.Range("O3").Value = 0.7967262871802238
$ws.Range("P3").Value = 0.7967262871802239
$ws.Range("Q3").Value = 311.89368457565
$ws.Range("R3").Value = 2807.04316118085
$ws.Range("S3").Value = 0.008939315414588787
$ws.Range("T3").Value = 0.008939315414588788
$ws.Range("G4").Value = 13.404345
$ws.Range("H4").Value = 40.213035
$ws.Range("I4").Value = 0.01122005832922476
$ws.Range("J4").Value = 0.01122005832922476
$ws.Range("O4").Value = 0.1365086454987356
$ws.Range("P4").Value = 0.1365086454987356
$ws.Range("Q4").Value = 53.43891008255499
$ws.Range("R4").Value = 480.9501907429949
$ws.Range("S4").Value = 0.001531634964939278
$ws.Range("T4").Value = 0.001531634964939278
$ws.Range("I5").Value = 0.9315566574535661
$ws.Range("J5").Value = 0.9315566574535661
$ws.Range("M5").Value = 1.949849666666667
$ws.Range("N5").Value = 5.849549000000001
$ws.Range("O5").Value = 0.06676506732104066
$ws.Range("P5").Value = 0.06676506732104066
$ws.Range("Q5").Value = 2170.005751591097
$ws.Range("R5").Value = 19530.05176431987
$ws.Range("S5").Value = 0.06219544294825095
$ws.Range("T5").Value = 0.06219544294825095
$ws.Range("I6").Value = 0.9315566574535661
$ws.Range("J6").Value = 0.9315566574535661
$ws.Range("O6").Value = 0.7967262871802238
$ws.Range("P6").Value = 0.7967262871802239
$ws.Range("S6").Value = 0.7421956769909992
$ws.Range("T6").Value = 0.7421956769909993
$ws.Range("I7").Value = 0.9315566574535661
$ws.Range("J7").Value = 0.9315566574535661
$ws.Range("O7").Value = 0.1365086454987356
$ws.Range("P7").Value = 0.1365086454987356
$ws.Range("S7").Value = 0.1271655375143159
$ws.Range("T7").Value = 0.1271655375143159
$ws.Range("I8").Value = 0.05722328421720919
$ws.Range("J8").Value = 0.05722328421720919
$ws.Range("M8").Value = 1.949849666666667
$ws.Range("N8").Value = 5.849549000000001
$ws.Range("O8").Value = 0.06676506732104066
$ws.Range("P8").Value = 0.06676506732104066
$ws.Range("Q8").Value = 133.2982324614706
$ws.Range("R8").Value = 1199.684092153235
$ws.Range("S8").Value = 0.003820516423093015
$ws.Range("T8").Value = 0.003820516423093015
$ws.Range("I9").Value = 0.05722328421720919
$ws.Range("J9").Value = 0.05722328421720919
$ws.Range("O9").Value = 0.7967262871802238
$ws.Range("P9").Value = 0.7967262871802239
$ws.Range("S9").Value = 0.04559129477463577
$ws.Range("T9").Value = 0.04559129477463578
$ws.Range("I10").Value = 0.05722328421720919
$ws.Range("J10").Value = 0.05722328421720919
$ws.Range("O10").Value = 0.1365086454987356
$ws.Range("P10").Value = 0.1365086454987356
$ws.Range("S10").Value = 0.007811473019480402
$ws.Range("T10").Value = 0.007811473019480402
